# Updated capital structure database
# Applies updated financial metrics to rows 2 and 3 of the earnings_debt sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that are identical for row 2 and row 3 (columns D through AQ),
# except for the cash_returned family (M,N,O,P,Q,R) which end up as -0 on row 3.
$commonValues = @{
    "D" = 0.0582
    "E" = 0.16
    "F" = -0.027
    "G" = 0.08927875243664718
    "H" = 0.08927875243664718
    "I" = 0.1072124756335283
    "J" = 0.09197084787031706
    "K" = 70.40000000000001
    "L" = 0.09148797920727746
    "U" = 157.1
    "V" = 0.4479612204163103
    "W" = 0.1714981729598051
    "X" = 0.100231193209154
    "Y" = 0.07126697975065116
    "Z" = 2.655096266648265
    "AA" = 0.2441914548209543
    "AB" = 0.08313903398391753
    "AC" = 0.1610524208370368
    "AD" = 101.6
    "AF" = 101.6
    "AG" = -55.5
    "AH" = 0.2246296705726288
    "AI" = 0.164854778516956
    "AJ" = -0.1880081300813008
    "AK" = -0.1208623693379791
    "AL" = 3.4
    "AM" = 3.4
    "AN" = 1.132664437012263
    "AO" = 24.26470588235294
    "AP" = -0.6187290969899666
    "AQ" = 24.26470588235294
}

foreach ($col in $commonValues.Keys) {
    $val = $commonValues[$col]
    $ws.Range("${col}2").Value = $val
    $ws.Range("${col}3").Value = $val
}

# cash_returned family: row 2 becomes 0, row 3 becomes -0 (numerically equal to 0)
foreach ($col in @("M", "N", "O", "P", "Q", "R")) {
    $ws.Range("${col}2").Value = 0
    $ws.Range("${col}3").Value = -0
}

# buybacks_cash_returned column (T) is removed entirely on both rows
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()
